# Updated symbol list on Wed Jan 25 21:49:13 UTC 2023 with GitHub Actions
# Refresh coin price / 1h-volume figures, and re-sync the coin list (rows 8-17)
# to the latest coinranking.com ordering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold text-formatted numbers/percentages
# (t="inlineStr" in the source). Force text number format before assignment so
# Excel doesn't silently re-type them as numeric/percentage values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '304.76'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.95%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '36.18'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.02%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.073'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.82%'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.30%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.944'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-3.40%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.832'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.70%'

$ws.Range("B8").Value = 'GateToken'

$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.063'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.64%'

$ws.Range("B9").Value = 'MXToken'

$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9333'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.68%'

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1517'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '39.54%'

$ws.Range("B11").Value = 'WazirX'

$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1911'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.18%'

$ws.Range("B12").Value = 'MandalaExchangeToken'

$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09016'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-6.23%'

$ws.Range("B13").Value = 'BitrueCoin'

$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03472'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-3.49%'

$ws.Range("B14").Value = 'BitMartToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09876'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.42%'

$ws.Range("B15").Value = 'BitForexToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001430'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.28%'

$ws.Range("B16").Value = 'TigerCash'

$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005795'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.36%'

$ws.Range("B17").Value = 'LEO'

$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.539'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.81%'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '3.09%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3444'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.66%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1282'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.52%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.025'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.44%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2388'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '8.89%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04471'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.60%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001204'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.78%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004867'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.96%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001224'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-2.19%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004418'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.77%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01990'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.28%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04821'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.13%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01046'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '8.83%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007340'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-6.28%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1368'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.04%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002063'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-2.55%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01068'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-7.69%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006089'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-6.03%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.62%'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '0.44%'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001184'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-9.00%'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.62%'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001990'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.62%'
